$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "pt_max" column (column F) entirely, shifting every
# column to its right one place to the left (G->F, H->G, ... M->L).
$ws.Columns("F:F").Delete()

# Reflect the author's post-edit selection: the whole (new) "boson" column F.
$ws.Range("F1:F1048576").Select()
